$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 2 gets a new country ("HU") and its year is set back to 2011.
$ws.Range("A2").Value = "HU"

# "2011" looks numeric, so set the cell to Text first so Excel stores it
# as a shared string (matching the source data) instead of a number, then
# drop the temporary number format so no extra style is left applied.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2011"
$ws.Range("B2").ClearFormats()
